$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the activity note text for row 14 (Task 16/18 entry)
$ws.Range("F14").Value = "Task 16: Complete. Task 18: Complete. Task 27: Preliminary work. Could never get it working. 30 minute meeting with Luca to resolve a merge conflict."

# Update hours worked for row 14 from 4 to 5
$ws.Range("E14").Value = 5

# Update the view: scroll position and active cell selection
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("F14").Select()
